$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "STEMMED" (no head) Proximity Search row (row 18) to be the
#    "Unmodified queries" variant with its updated figures.
$ws.Range("A18").Value = "Proximity Search (Unmodified queries)"
$ws.Range("B18").Value = 0.2354
$ws.Range("C18").Value = 0.36
$ws.Range("D18").Value = 0.2907

# 2. Remove the now-redundant standalone "Proximity Search (Unmodified queries)"
#    row that used to sit under the "With HEAD" STEMMED block (old row 25).
#    Deleting it shifts everything below up by one row.
$ws.Rows(25).Delete()

# 3. The "NON STEMMED" (no head) Proximity Search row -- originally row 32,
#    now row 31 after the deletion above -- becomes the "Unmodified queries"
#    variant with its updated figures.
$ws.Range("A31").Value = "Proximity Search (Unmodified queries)"
$ws.Range("B31").Value = 0.1842
$ws.Range("C31").Value = 0.272
$ws.Range("D31").Value = 0.2347

# 4. Refresh the saved view state to reflect where the user ended up editing.
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("A31").Select()
